{"js": "// Thesis writing plan edits:\n// 1) Collapse split/proofed runs into single clean runs (no text change).\n// 2) Insert a new \"Acknowledgements and statement of contribution\" list\n//    item right after \"Write abstract (half day)\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst fixups = [\n  {\n    match: \"Restructure body chapters into\",\n    text: \"Restructure body chapters into cosmics data, x-ray data, comparison and edit (2 days)\"\n  },\n  {\n    match: \"Body formatting,\",\n    text: \"Body formatting, eg. positioning figures (half day)\"\n  },\n  {\n    match: \"Brigitte full edit feedback\",\n    text: \"Brigitte full edit feedback ( 3 days)\"\n  }\n];\n\nlet abstractPara = null;\n\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const text = para.text;\n\n  const fixup = fixups.find((f) => text.indexOf(f.match) === 0);\n  if (fixup) {\n    para.insertText(fixup.text, \"Replace\");\n  }\n\n  if (text.indexOf(\"Write abstract\") === 0) {\n    abstractPara = para;\n  }\n}\n\nawait context.sync();\n\nif (abstractPara) {\n  abstractPara.insertParagraph(\n    \"Acknowledgements and statement of contribution (few hours)\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "# Thesis writing plan edits:\n# 1) Collapse split/proofed runs into single clean runs (no text change).\n# 2) Insert a new \"Acknowledgements and statement of contribution\" list\n#    item right after \"Write abstract (half day)\".\n\n$d = $word.ActiveDocument\n\nfunction Set-CleanParagraphText([int]$paraIndex, [string]$newText) {\n    $paragraph = $d.Paragraphs.Item($paraIndex)\n    $start = $paragraph.Range.Start\n    $end = $paragraph.Range.End - 1\n    $r = $d.Range($start, $end)\n    $r.Delete()\n    $p2 = $d.Paragraphs.Item($paraIndex)\n    $ip = $p2.Range\n    $ip.Collapse(1)\n    $ip.InsertBefore($newText)\n}\n\n$fixups = @(\n    @{ Match = \"Restructure body chapters into\"; Text = \"Restructure body chapters into cosmics data, x-ray data, comparison and edit (2 days)\" },\n    @{ Match = \"Body formatting,\"; Text = \"Body formatting, eg. positioning figures (half day)\" },\n    @{ Match = \"Brigitte full edit feedback\"; Text = \"Brigitte full edit feedback ( 3 days)\" }\n)\n\n$count = $d.Paragraphs.Count\n$abstractIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n\n    foreach ($fixup in $fixups) {\n        if ($t.StartsWith($fixup.Match)) {\n            Set-CleanParagraphText $i $fixup.Text\n        }\n    }\n\n    if ($t.StartsWith(\"Write abstract\")) {\n        $abstractIndex = $i\n    }\n}\n\nif ($abstractIndex -ge 1) {\n    $abstractPara = $d.Paragraphs.Item($abstractIndex)\n    $abstractPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($abstractIndex + 1)\n    $newPara.Range.Text = \"Acknowledgements and statement of contribution (few hours)\"\n}\n"}
